# Update the Key Outbreaks cluster table: the source data refreshed with
# new clusters added, resolved clusters removed, and active-case counts
# updated for the remaining clusters.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('3398 BlueCross Elly Kay Mordialloc', 32),
    @('3564 Waverley Valley Aged Care Glen Waverley', 10),
    @('3601 Baptcare Westhaven community', 14),
    @('3646 Mornington Bay Care Community Mount Martha', 14),
    @('3653 Fronditha Thalpori St Albans Aged Care', 21),
    @('3975 Aurrum Aged Care Brunswick West', 13),
    @('4257 BlueCross The Gables Camberwell', 18),
    @('4295 Hope Aged Care Sunshine West', 18),
    @('4314 Estia Health Ardeer', 18),
    @('44304 Brighton Primary School Brighton', 15),
    @('44380 Plenty Parklands Primary School Mill Park', 10),
    @('44414 Buangor Primary School Buangor', 20),
    @('44584 Badger Creek Primary School Badger Creek', 13),
    @('44593 Torquay P-6 College Torquay', 17),
    @('44666 Gardenvale Primary School Senior School Campus Brighton East', 10),
    @('44761 Coburg North Primary School Coburg', 22),
    @('44828 Cheltenham East Primary School Cheltenham', 10),
    @('44930 St Albans Heights Primary School', 14),
    @('44979 Campbellfield Heights Primary School Campbellfield', 18),
    @('45248 Brookside P-9 College Caroline Springs', 15),
    @('45265 Cranbourne Carlisle Primary School', 10),
    @('45350 Ballarat High School Lake Gardens', 10),
    @('45648 St Brendans Primary School Shepparton', 11),
    @('45677 St. Peter''s Primary School Bentleigh East', 10),
    @('45708 St Augustine''s Primary School Yarraville', 10),
    @('45719 St Joseph''s Primary School Numurkah', 12),
    @('4574 Village Glen Aged Care Residences Mornington', 10),
    @('45902 Mother of God Primary School Ardeer', 12),
    @('45988 St Macartan''s Parish Primary School', 11),
    @('46052 St. Francis of Assisi Primary School Mill Park', 13),
    @('46115 St Luke''s Catholic Primary School Shepparton North', 11),
    @('46116 Saint Francis Of Assisi Primary School Baranduda', 10),
    @('46277 St Paul''s Anglican Grammar School Warragul', 11),
    @('46306 King''s College Warrnambool', 10),
    @('46328 Ilim College Inverloch Crescent Dallas', 12),
    @('50722 Chairo Christian School Leongatha Campus Leongatha', 11),
    @('Kororoit Creek Primary School Burnside Heights', 20),
    @('Mambourin Enterprises Allara Deer Park', 20),
    @('North St Kilda Children''s Centre St Kilda', 10),
    @('Oakleigh South Primary School Oakleigh South', 12),
    @('Rosebud Primary School Rosebud', 16),
    @('Springside Primary School Caroline Springs', 13),
    @('St Christophers Primary School Airport West', 12),
    @('St Louis de Montfort''s School Aspendale', 16),
    @('Torquay Hotel Torquay', 28),
)

# Clear out the old table body (below the header row) before writing
# back the refreshed row count.
$oldUsed = $ws.UsedRange
$lastRow = $oldUsed.Row + $oldUsed.Rows.Count - 1
if ($lastRow -ge 2) {
    $ws.Range("A2:B$lastRow").Clear()
}

$rowCount = $data.Count
for ($i = 0; $i -lt $rowCount; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $data[$i][0]
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
}
